# Fixed Bento 80 Test scripts
# Appends "order By ... LIMIT 100" clauses to the three Cypher/SQL query
# strings shown on the "startup" sheet (CasesTab / SamplesTab / FilesTab rows),
# and updates the two affected rows' heights so the wrapped text keeps fitting.
#
# NOTE: cells are updated in B2, B3, B4 order (top -> bottom) so that the
# workbook's shared-string table regenerates with the same relative ordering
# as the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = "MATCH (ss:study_subject)`nMATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`nWITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`nMATCH (ss)<-[:sf_of_study_subject]-(sf)`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`nMATCH (ss)<-[:demographic_of_study_subject]-(demo)`nWHERE tp.chemotherapy_regimen in [`"Taxane only`"]`nreturn ss.study_subject_id as ``Case ID``,`n       p.program_acronym as ``Program Code``,`n        p.program_id as Program_ID,`n       s.study_acronym as ``Arm``,`n       ss.disease_subtype as ``Diagnosis``,`n       sf.grouped_recurrence_score AS ``Recurrence Score``,`n       d.tumor_size_group AS ``tumor_size``,`n       d.er_status AS ``ER Status``,`n       d.pr_status AS ``PR Status``,`n       demo.age_at_index AS ``Age (years)``,`ndemo.survival_time AS ``Survival (days)```n order By ss.study_subject_id ASC LIMIT 100 "

$samplesQuery = "MATCH (ss:study_subject)`nWITH COLLECT(ss.study_subject_id) AS all_subjects`nMATCH (samp:sample)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`nMATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`nWHERE tp.chemotherapy_regimen IN [`"Taxane only`"]`nWITH`n    distinct lp,`n    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,`n    collect(distinct f.file_id) AS files,`n    samp, ss, s, p, all_subjects`nRETURN`n samp.sample_id AS ``Sample ID``,`n            ss.study_subject_id AS ``Case ID``,`n            p.program_acronym AS ``Program Code``,`n            s.study_acronym AS ``Arm``,`n            ss.disease_subtype AS ``Diagnosis``,`n            samp.tissue_type AS ``Tissue Type``,`n            samp.composition AS ``Tissue Composition``,`n            samp.sample_anatomic_site AS ``Sample Anatomic Site``,`n            samp.method_of_sample_procurement AS ``Sample Procurement Method```n order By samp.sample_id ASC LIMIT 100"

$filesQuery = "MATCH (f:file)-->(parent)`nMATCH (f)-[:file_of_sample]->(samp)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)`nMATCH (s)-[:study_of_program]->(p)`nMATCH (d)-[:diagnosis_of_study_subject]->(ss)`nMATCH (tp)-[:tp_of_diagnosis]->(d)`nWHERE tp.chemotherapy_regimen IN [`"Taxane only`"]`nWITH`n        f, parent,p, ss, d,tp, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent,p, ss, d,tp, s, samp,`n        f.file_size /(1024^i) AS value,`n        10^precision AS factor,`n        units[i] as unit`nWITH`n        f, parent,p, ss, d,tp, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN Distinct`n    f.file_name AS ``File Name``,`n    head(labels(samp)) AS ``Association``,`n    f.file_description AS ``Description``,`n    f.file_format AS ``File Format``,`n     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    p.program_acronym AS ``Program Code``,`n    s.study_acronym AS ``Arm``,`n    ss.study_subject_id AS ``Case ID``,`n    samp.sample_id AS ``Sample ID```n  order By f.file_name ASC LIMIT 100"

$ws.Cells.Item(2, 2).Value2 = $casesQuery
$ws.Cells.Item(3, 2).Value2 = $samplesQuery
$ws.Cells.Item(4, 2).Value2 = $filesQuery

# Row heights grew slightly to accommodate the extra wrapped line of text
# (row 4 was already at Excel's 409.6pt wrap cap, so it is left untouched).
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360
